$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tables")

# ---------------------------------------------------------------------------
# 1. Data correction: ids 1 & 2 ("bob", "jon") change role from operator ->
#    manager in the dataTable.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "manager"
$ws.Range("C3").Value = "manager"

# ---------------------------------------------------------------------------
# 2. New workbook-level defined names used by the new examples.
# ---------------------------------------------------------------------------
$wb.Names.Add("role_in", '=tables!$C$26')
$wb.Names.Add("subrange", '=INDEX(dataTable[salary], MATCH(role_in, dataTable[role], 0)): INDEX(dataTable[salary], MATCH(role_in, dataTable[role], 1))')

# ---------------------------------------------------------------------------
# 3. Move the existing "array formula multiplication" example down from
#    rows 25:29 to rows 38:42 to make room for the new examples.
# ---------------------------------------------------------------------------
$ws.Range("D25:D29").Copy()
$ws.Range("D38:D42").PasteSpecial(-4122)
$eNote = $ws.Range("E25").Value2

$ws.Range("D25").CurrentArray.ClearContents()
$ws.Range("E25").ClearContents()
# Touch the old array's anchor cell first so the rest of its former range is
# released for plain writes (engine quirk: cells that were part of an array
# formula stay locked for direct writes until the anchor itself is written).
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0

$ws.Range("D38:D42").FormulaArray = '=dataTable[age] * dataTable[salary]'
$ws.Range("E38").Value = $eNote

# ---------------------------------------------------------------------------
# 4. New example block: "SELECT/SUM WHERE using array formula" (rows 25-36).
#    Values are written in the exact order the original authoring produced
#    them so new shared-string entries land at the same indices.
# ---------------------------------------------------------------------------

# Row 25 - headers (bold, no border) styled like G1/H1.
$ws.Range("G1").Copy()
$ws.Range("C25:D25").PasteSpecial(-4122)

# Row 26 - the role_in input cell (styled like C18/C21) + SUM(subrange).
$ws.Range("C18").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("D26").PasteSpecial(-4122)

# Row 27 cell formatting
$ws.Range("A2").Copy()
$ws.Range("D27").PasteSpecial(-4122)

# Row 28 placeholder formatting (C28:C36) based on C18 minus the border.
$ws.Range("C18").Copy()
$ws.Range("C28:C36").PasteSpecial(-4122)
$ws.Range("C28:C36").Borders.LineStyle = -4142

$ws.Range("A2").Copy()
$ws.Range("D28:D32").PasteSpecial(-4122)

# --- values, in shared-string creation order: C25, E26, D25, E28, E27 ---
$ws.Range("C25").Value = "role_in"
$ws.Range("E26").Value = "(INDEX, MATCH) : (INDEX, MATCH)"
$ws.Range("D25").Value = "subrange"
$ws.Range("E28").Value = "SELECT WHERE using array formula"
$ws.Range("E27").Value = "SUM IF WHERE using array formula"

# --- formulas / remaining values ---
$ws.Range("C26").Value = "manager"
$ws.Range("D26").Formula = '=SUM(subrange)'
$ws.Range("D27").FormulaArray = '=SUM(IF(dataTable[role]=role_in, dataTable[salary], 0))'
$ws.Range("C28:C36").ClearContents()
$ws.Range("D28:D32").FormulaArray = '=IF(dataTable[role]=role_in, dataTable[salary], 0)'

# ---------------------------------------------------------------------------
# 5. View-state: make "tables" the active sheet / tab, and set its
#    selection to match the new layout.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C27").Select()

Write-Host "done"
